$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the existing "IP" header cell (H1) onto the two
# new header cells I1 and J1 so they pick up the same bold/border/centered
# style (s="1") used by the rest of the header row.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# New header labels
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for rows 2..59 (column I = "I0", column J = "IF")
$iVals = @(8,8,7,9,7,6,7,9,5,7,11,6,6,9,6,7,7,9,6,6,8,7,6,7,5,7,7,7,7,8,6,7,7,7,8,8,3,7,7,7,8,5,5,2,9,7,8,7,8,8,7,6,5,7,4,5,4,3)
$jVals = @(8,8,8,9,7,6,7,9,5,7,11,6,6,9,6,7,7,9,6,6,8,8,6,7,5,8,7,7,7,8,6,7,7,7,8,8,3,7,7,7,8,5,5,2,9,7,8,7,8,8,7,6,5,7,4,5,4,3)

for ($idx = 0; $idx -lt $iVals.Length; $idx++) {
    $row = $idx + 2
    $ws.Cells.Item($row, 9).Value = $iVals[$idx]
    $ws.Cells.Item($row, 10).Value = $jVals[$idx]
}
